$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the first sheet
$ws.Name = "Export as TSV"

# Freeze the header row (row 1)
$ws.Activate()
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true

# Add errorTitle/error messages to the existing data validations
$validations = @(
    @{ Range = "I2:I1048576"; Title = "Value must come from list"; Message = "Value must be one of: mass_spectrometry_imaging." },
    @{ Range = "J2:J1048576"; Title = "Value must come from list"; Message = "Value must be one of: MALDI-IMS." },
    @{ Range = "K2:K1048576"; Title = "Value must come from list"; Message = "Value must be one of: protein / metabolites / lipids." },
    @{ Range = "L2:L1048576"; Title = "Not a boolean"; Message = 'The values in this column must be "TRUE" or "FALSE".' },
    @{ Range = "O2:O1048576"; Title = "Value must come from list"; Message = "Value must be one of: MALDI / MALDI-2 / DESI / SIMS / nESI." },
    @{ Range = "P2:P1048576"; Title = "Value must come from list"; Message = "Value must be one of: negative ion mode / positive ion mode." },
    @{ Range = "Q2:Q1048576"; Title = "Not a number"; Message = "The values in this column must be numbers." },
    @{ Range = "R2:R1048576"; Title = "Not a number"; Message = "The values in this column must be numbers." },
    @{ Range = "S2:S1048576"; Title = "Not a number"; Message = "The values in this column must be numbers." },
    @{ Range = "T2:T1048576"; Title = "Value must come from list"; Message = "Value must be one of: nm / um." },
    @{ Range = "U2:U1048576"; Title = "Not a number"; Message = "The values in this column must be numbers." },
    @{ Range = "V2:V1048576"; Title = "Value must come from list"; Message = "Value must be one of: nm / um." }
)

foreach ($v in $validations) {
    $rng = $ws.Range($v.Range)
    $rng.Validation.ErrorTitle = $v.Title
    $rng.Validation.ErrorMessage = $v.Message
}
